# Fruta / hortaliza, semanal
# A new weekly price-report row is inserted at row 65 (pushing the
# existing rows 65-141 down to 66-142); the new row carries a fresh
# Frutilla price observation for "Provincia de Melipilla".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 65, shifting rows
# 65..141 down to 66..142 (dimension/used-range updates automatically).
$ws.Rows(65).Insert()

# Populate the newly-inserted row 65 with the new observation.
$ws.Cells.Item(65, 1).Value  = 10
$ws.Cells.Item(65, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(65, 3).Value  = "La Araucanía"
$ws.Cells.Item(65, 4).Value  = 44482
$ws.Cells.Item(65, 5).Value  = 9
$ws.Cells.Item(65, 6).Value  = "Fruta"
$ws.Cells.Item(65, 7).Value  = 100101
$ws.Cells.Item(65, 8).Value  = "Berries"
$ws.Cells.Item(65, 9).Value  = 100112025
$ws.Cells.Item(65, 10).Value = "Frutilla"
$ws.Cells.Item(65, 11).Value = "Sin especificar"
$ws.Cells.Item(65, 12).Value = "Primera"
$ws.Cells.Item(65, 13).Value = 110
$ws.Cells.Item(65, 14).Value = 12000
$ws.Cells.Item(65, 15).Value = 12000
$ws.Cells.Item(65, 16).Value = 12000
$ws.Cells.Item(65, 17).Value = "`$/bandeja 7 kilos"
$ws.Cells.Item(65, 18).Value = "Provincia de Melipilla"
$ws.Cells.Item(65, 19).Value = 1714
$ws.Cells.Item(65, 20).Value = 7
